# Add "5V 4A DC/DC Power Supply" (Recon Power R-745.0P) as a new BOM line
# in row 11, which previously was a blank spacer row between the
# connector/IC section and the PCB/resistor section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Part number / manufacturer / supplier / description
$ws.Range("A11").Value = "R-745.0P"
$ws.Range("B11").Value = "Recon Power"
$ws.Range("C11").Value = "Digikey"
$ws.Range("D11").Value = "5V 4A Output Power Supply"

# Quantity / unit price
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 35.69

# Extended price = Quantity * Price (same pattern used by the rest of the
# G column in this block)
$ws.Range("G11").Formula = "=E11*F11"

# Give the new row's "status" cell (column H) the same yellow highlight
# used for other newly-added / not-yet-received parts, by recoloring the
# fill while keeping the existing status-column font.
$ws.Range("H11").Interior.Color = 65535
